# Update need_to_buy sheet data: rows 2-15, columns A-F
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45909, 5903.89122679882, 5148.10730071727, 5688, 6508.829071, 2.71021437160205),
    @(45910, 5914.10694019361, 5178.09697868103, 2952, 6521.923268, 118.079721103643),
    @(45911, 5913.30581487482, 5171.56408168856, 2952, 6521.069439, 117.805321075573),
    @(45912, 5913.14141862252, 4470.46864981824, 2952, 6520.667157, 88.5830995081547),
    @(45913, 1177.4817126383, 2215.40132970979, 2952, 2248.970378, 13.9537497946452),
    @(45914, 1023.06998026009, 2162.06663960271, 2952, 2036.074808, 9.29464447260932),
    @(45915, 6061.80298841384, 5033.79850517088, 2952, 6722.443532, 114.26829369821),
    @(45916, 6061.80298841384, 5232.79287965655, 2952, 6722.443532, 122.559725968446),
    @(45917, 6061.80298841384, 5246.36320477272, 2952, 6722.443532, 123.12515618162),
    @(45918, 6061.80298841384, 5398.25897984135, 2952, 6722.443532, 129.454146809479),
    @(45919, 6061.80298841384, 4951.53176359997, 2952, 6722.443532, 110.840512799422),
    @(45920, 1187.6944176276, 2907.8570871323, 2952, 2248.908814, 42.3779784793627),
    @(45921, 1056.51488323291, 2720.24486780317, 2952, 2106.592241, 34.0967593987609),
    @(45922, 5979.99309220113, 5084.23316341545, 2952, 6713.304478, 119.397689550597)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $row++
}

$wb.Save()
